# Apply the "cryptos" data refresh described by the commit:
#   "Updated cryptos list on Sat Jun  3 22:24:45 UTC 2023 with GitHub Actions"
# Column D = Price, Column E = Volume(1h). All cells on this sheet are plain
# text (several "prices" use dotted thousands-separators like "27.112.75" and
# are not valid numbers at all), so every write below is forced to Text first
# and the style is reset back to Normal afterwards so no stray number format
# is left attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $ref, $text) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Column D: Price ---
Set-TextValue $ws "D2" "27.112.75"
Set-TextValue $ws "D3" "1.894.53"
Set-TextValue $ws "D5" "306.70"
Set-TextValue $ws "D7" "0.5209"
Set-TextValue $ws "D8" "0.3761"
Set-TextValue $ws "D9" "0.07266"
Set-TextValue $ws "D10" "21.14"
Set-TextValue $ws "D11" "0.8988"
Set-TextValue $ws "D12" "0.08193"
Set-TextValue $ws "D13" "1.957.31"
Set-TextValue $ws "D14" "96.23"
Set-TextValue $ws "D15" "5.297"
Set-TextValue $ws "D17" "0.000008601"
Set-TextValue $ws "D18" "14.57"
Set-TextValue $ws "D19" "1.003"
Set-TextValue $ws "D20" "27.137.39"
Set-TextValue $ws "D21" "5.086"
Set-TextValue $ws "D22" "10.70"
Set-TextValue $ws "D23" "6.410"
Set-TextValue $ws "D24" "148.48"
Set-TextValue $ws "D25" "2.305"
Set-TextValue $ws "D26" "18.19"
Set-TextValue $ws "D27" "1.739"
Set-TextValue $ws "D28" "115.14"
Set-TextValue $ws "D29" "4.790"
Set-TextValue $ws "D30" "4.865"
Set-TextValue $ws "D31" "0.09204"
Set-TextValue $ws "D32" "0.05035"
Set-TextValue $ws "D33" "0.7912"
Set-TextValue $ws "D34" "1.216"
Set-TextValue $ws "D35" "3.427"
Set-TextValue $ws "D36" "2.976"
Set-TextValue $ws "D37" "2.610"
Set-TextValue $ws "D38" "0.5719"
Set-TextValue $ws "D39" "0.01991"
Set-TextValue $ws "D41" "9.029"
Set-TextValue $ws "D42" "6.549"
Set-TextValue $ws "D43" "116.42"
Set-TextValue $ws "D44" "0.1512"
Set-TextValue $ws "D45" "0.4854"
Set-TextValue $ws "D47" "10.09"
Set-TextValue $ws "D48" "1.620"
Set-TextValue $ws "D49" "38.19"
Set-TextValue $ws "D50" "63.58"
Set-TextValue $ws "D51" "0.05929"

# --- Column E: Volume(1h) (values keep their original 2-space padding) ---
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -2.27%  "
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  -0.27%  "
